# The deck currently applies the "Integral" (Red Violet) design to the
# slides, while the notes master still carries the old "Office Theme"
# colour scheme. The author flipped which theme drives the main
# presentation: the slides (and slide master) should now use the plain
# "Office" colour palette.
#
# Re-colour the active theme (the one the slide master / presentation
# actually uses) to the stock Office theme colours.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # Dark 1      -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # Light 1     -> FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # Dark 2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # Light 2     -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # Accent 1    -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # Accent 2    -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # Accent 3    -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # Accent 4    -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # Accent 5    -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # Accent 6    -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # Hyperlink   -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # Followed Hyperlink -> 954F72
